# Edit script: populate rainfall values in column B and apply top-vertical alignment style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row (2-366) -> new rainfall_mm_m2 value for column B.
# All of these rows also receive "vertical top" alignment (matching the new style added to cellXfs).
$data = @(
    @(2, 1.437),
    @(3, 2.458),
    @(4, 5.783),
    @(5, 8.896),
    @(6, 18.665),
    @(7, 5.241),
    @(8, 7.258),
    @(9, 10.287),
    @(10, 5.73),
    @(11, 5.72),
    @(12, 10.219),
    @(13, 1.323),
    @(14, 2.383),
    @(15, 4.869),
    @(16, 4.281),
    @(17, 8.446),
    @(18, 9.064),
    @(19, 7.844),
    @(20, 18.157),
    @(21, 7.414),
    @(22, 11.636),
    @(23, 13.683),
    @(24, 3.201),
    @(25, 6.972),
    @(26, 11.562),
    @(27, 6.561),
    @(28, 15.327),
    @(29, 2.676),
    @(30, 8.652),
    @(31, 24.894),
    @(32, 30.123),
    @(33, 19.222),
    @(34, 3.631),
    @(35, 32.277),
    @(36, 9.034),
    @(37, 4.512),
    @(38, 18.005),
    @(39, 10.156),
    @(40, 6.494),
    @(41, 13.51),
    @(42, 11.583),
    @(43, 1.205),
    @(44, 6.548),
    @(45, 14.007),
    @(46, 9.626),
    @(47, 2.126),
    @(48, 11.055),
    @(49, 20.735),
    @(50, 17.884),
    @(51, 11.419),
    @(52, 4.341),
    @(53, 7.55),
    @(54, 6.809),
    @(55, 12.474),
    @(56, 11.781),
    @(57, 6.74),
    @(58, 22.599),
    @(59, 3.813),
    @(60, 15.529),
    @(61, 6.843),
    @(62, 19.179),
    @(63, 4.347),
    @(64, 9.799),
    @(65, 7.884),
    @(66, 26.105),
    @(67, 3.626),
    @(68, 4.844),
    @(69, 26.005),
    @(70, 25.885),
    @(71, 15.754),
    @(72, 10.73),
    @(73, 0.452),
    @(74, 11.654),
    @(75, 3.607),
    @(76, 24.537),
    @(77, 0.989),
    @(78, 9.869),
    @(79, 2.626),
    @(80, 3.913),
    @(81, 0.411),
    @(82, 0),
    @(83, 0.228),
    @(84, 0),
    @(85, 9.76),
    @(86, 0.022),
    @(87, 8.069),
    @(88, 0.146),
    @(89, 5.603),
    @(90, 4.576),
    @(91, 0.923),
    @(92, 0),
    @(93, 0.003),
    @(94, 0.241),
    @(95, 13.476),
    @(96, 10.672),
    @(97, 5.079),
    @(98, 5.899),
    @(99, 16.653),
    @(100, 0.859),
    @(101, 13.634),
    @(102, 7.866),
    @(103, 4.024),
    @(104, 0.002),
    @(105, 0.016),
    @(106, 0),
    @(107, 1.154),
    @(108, 0.347),
    @(109, 0),
    @(110, 0),
    @(111, 0),
    @(112, 0),
    @(113, 0.036),
    @(114, 0),
    @(115, 0.803),
    @(116, 0.004),
    @(117, 0),
    @(118, 0),
    @(119, 0),
    @(120, 0),
    @(121, 0),
    @(122, 0),
    @(123, 0),
    @(124, 0),
    @(125, 0),
    @(126, 0),
    @(127, 0),
    @(128, 0),
    @(129, 0),
    @(130, 0),
    @(131, 0),
    @(132, 0),
    @(133, 0),
    @(134, 0),
    @(135, 1.403),
    @(136, 0),
    @(137, 0.001),
    @(138, 0.057),
    @(139, 0.031),
    @(140, 0.392),
    @(141, 10.333),
    @(142, 0.28),
    @(143, 0.001),
    @(144, 0),
    @(145, 0.001),
    @(146, 0),
    @(147, 0),
    @(148, 0),
    @(149, 0.001),
    @(150, 0),
    @(151, 0),
    @(152, 0),
    @(153, 0),
    @(154, 0),
    @(155, 0),
    @(156, 0),
    @(157, 0),
    @(158, 0),
    @(159, 0),
    @(160, 0),
    @(161, 0),
    @(162, 1.136),
    @(163, 0),
    @(164, 0),
    @(165, 0),
    @(166, 0),
    @(167, 0),
    @(168, 0),
    @(169, 0),
    @(170, 0),
    @(171, 0),
    @(172, 0),
    @(173, 0),
    @(174, 0),
    @(175, 0),
    @(176, 0),
    @(177, 0),
    @(178, 0),
    @(179, 0),
    @(180, 0),
    @(181, 0),
    @(182, 0),
    @(183, 0),
    @(184, 0),
    @(185, 0.006),
    @(186, 0),
    @(187, 0),
    @(188, 0),
    @(189, 0),
    @(190, 0),
    @(191, 0),
    @(192, 0),
    @(193, 0),
    @(194, 0),
    @(195, 0),
    @(196, 0),
    @(197, 0),
    @(198, 0),
    @(199, 0),
    @(200, 0),
    @(201, 0),
    @(202, 0),
    @(203, 0),
    @(204, 0),
    @(205, 0),
    @(206, 0),
    @(207, 0),
    @(208, 0),
    @(209, 0),
    @(210, 0),
    @(211, 0),
    @(212, 0),
    @(213, 0),
    @(214, 0),
    @(215, 0),
    @(216, 0),
    @(217, 0),
    @(218, 0),
    @(219, 0),
    @(220, 0),
    @(221, 0),
    @(222, 0),
    @(223, 0),
    @(224, 0),
    @(225, 0),
    @(226, 0.001),
    @(227, 0),
    @(228, 0),
    @(229, 0),
    @(230, 0.017),
    @(231, 0),
    @(232, 0),
    @(233, 0),
    @(234, 0.01),
    @(235, 0.284),
    @(236, 0),
    @(237, 0),
    @(238, 0),
    @(239, 0),
    @(240, 0.004),
    @(241, 0),
    @(242, 0),
    @(243, 0),
    @(244, 0),
    @(245, 0),
    @(246, 0),
    @(247, 0),
    @(248, 0),
    @(249, 0),
    @(250, 0),
    @(251, 0),
    @(252, 0),
    @(253, 0),
    @(254, 0),
    @(255, 0),
    @(256, 0),
    @(257, 0),
    @(258, 0),
    @(259, 0),
    @(260, 0),
    @(261, 0.001),
    @(262, 0),
    @(263, 0.008),
    @(264, 0.017),
    @(265, 0.047),
    @(266, 0),
    @(267, 0),
    @(268, 0),
    @(269, 0),
    @(270, 0),
    @(271, 0),
    @(272, 0),
    @(273, 0),
    @(274, 0),
    @(275, 0),
    @(276, 0),
    @(277, 0),
    @(278, 0),
    @(279, 0.065),
    @(280, 0),
    @(281, 0),
    @(282, 0),
    @(283, 0),
    @(284, 0),
    @(285, 0),
    @(286, 0),
    @(287, 0),
    @(288, 0),
    @(289, 0),
    @(290, 0),
    @(291, 0),
    @(292, 0),
    @(293, 0),
    @(294, 0),
    @(295, 0),
    @(296, 0),
    @(297, 0),
    @(298, 0),
    @(299, 0.191),
    @(300, 0.325),
    @(301, 6.636),
    @(302, 1.354),
    @(303, 0),
    @(304, 0),
    @(305, 0),
    @(306, 0),
    @(307, 0),
    @(308, 4.035),
    @(309, 5.354),
    @(310, 22.962),
    @(311, 8.174),
    @(312, 15.566),
    @(313, 8.498),
    @(314, 6.694),
    @(315, 0.001),
    @(316, 0),
    @(317, 0.035),
    @(318, 0),
    @(319, 0),
    @(320, 0),
    @(321, 0),
    @(322, 0.005),
    @(323, 0),
    @(324, 0.016),
    @(325, 0),
    @(326, 1.847),
    @(327, 0),
    @(328, 0.021),
    @(329, 0.009),
    @(330, 0.497),
    @(331, 8.8),
    @(332, 12.929),
    @(333, 7.724),
    @(334, 1.808),
    @(335, 0.232),
    @(336, 0.045),
    @(337, 0.137),
    @(338, 0),
    @(339, 0),
    @(340, 0.019),
    @(341, 0),
    @(342, 0),
    @(343, 0.406),
    @(344, 0.985),
    @(345, 1.05),
    @(346, 1.712),
    @(347, 3.93),
    @(348, 4.319),
    @(349, 0.016),
    @(350, 8.542),
    @(351, 6.588),
    @(352, 0.128),
    @(353, 1.057),
    @(354, 10.374),
    @(355, 5.537),
    @(356, 8.413),
    @(357, 17.421),
    @(358, 14.303),
    @(359, 3.344),
    @(360, 5.979),
    @(361, 4.012),
    @(362, 1.105),
    @(363, 2.464),
    @(364, 2.131),
    @(365, 2.801),
    @(366, 1.903)
)

foreach ($item in $data) {
    $row = $item[0]
    $val = $item[1]
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $val
    $cell.VerticalAlignment = -4160
}

# Update the active selection / scroll position to match the saved view.
$ws.Range("E358").Select()
$excel.ActiveWindow.ScrollRow = 349
$excel.ActiveWindow.ScrollColumn = 1
